# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" / "Valor Mora" table on Hoja1 (rows 16-61, columns E:F)
# previously listed account-statement periods in ascending order
# (1910, 1911, ... 2307). The author removed the old account statements and
# appended the newest ones, which - once the underlying database/table is
# rebuilt - leaves the periods listed in descending order
# (2307, 2306, ... 1910) while the set of "Valor Mora" amounts stays the
# same (forty-five periods at 40000 and one period at 34666), just
# re-attached to the now-reordered period labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Original ascending period labels (row 16 .. row 61)
$periods = @(
    "1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
    "2301","2302","2303","2304","2305","2306","2307"
)

# Original "Valor Mora" amounts (row 16 .. row 61): every period is 40000
# except the most recent one (2307), which is 34666.
$values = @()
for ($i = 0; $i -lt $periods.Length; $i++) { $values += 40000 }
$values[$values.Length - 1] = 34666

# Reverse both lists so the newest period ends up first (row 16) and the
# oldest ends up last (row 61), carrying its own Valor Mora amount with it.
$lastIndex = $periods.Length - 1
$periodsDesc = $periods[$lastIndex..0]
$valuesDesc = $values[$lastIndex..0]

$firstRow = 16
for ($i = 0; $i -lt $periodsDesc.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $periodsDesc[$i]
    $ws.Range("F$row").Value = $valuesDesc[$i]
}
